$wb = $excel.ActiveWorkbook

# --- Rename the "Include from LOINC" worksheet to "Include #0" ---
$includeSheet = $wb.Worksheets.Item("Include from LOINC")
$includeSheet.Name = "Include #0"

# --- Update the Metadata sheet ---
$metaSheet = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$metaSheet.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"

# Insert a new row above row 11 ("Description") for the new "Jurisdiction" property,
# shifting Description/Purpose/Copyright/Immutable down by one row.
# Use CopyOrigin = xlFormatFromLeftOrAbove (0) so the new row inherits the
# formatting of the row above it instead of creating stray new styles.
$metaSheet.Rows.Item(11).Insert(-4121, 0)

$metaSheet.Cells.Item(11, 1).Value = "Jurisdiction"
$metaSheet.Cells.Item(11, 2).Value = ""

# Re-apply the exact formatting used by the other property rows (row style "s=2")
# by copying it, cell by cell, from the row that now holds "Description"
# (row 12) — whole-row copy/paste tends to introduce extra style entries,
# but single-cell format copy reuses the existing style cleanly.
$metaSheet.Cells.Item(12, 1).Copy()
$metaSheet.Cells.Item(11, 1).PasteSpecial(-4122)
$metaSheet.Cells.Item(12, 2).Copy()
$metaSheet.Cells.Item(11, 2).PasteSpecial(-4122)
